{"js": "// Commit: \"adicionando atributo de situacao para pedido - futuro enum\"\n// The \"Pedido (...)\" bullet line needs a new \", situacao\" column appended\n// right after \"id_carrinho\" (and before the closing parenthesis), matching\n// the plain (non-italic) formatting used for the other regular attributes\n// in that line (only \"id_carrinho\" itself stays italic, marking it as a\n// foreign key).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the \"- Pedido (...)\" definition paragraph explicitly instead of\n// assuming a fixed index, so the script is resilient to minor reflow.\nlet pedidoPara = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const text = paragraphs.items[i].text;\n  if (text.indexOf(\"Pedido (\") !== -1 && text.indexOf(\"id_carrinho\") !== -1) {\n    pedidoPara = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!pedidoPara) {\n  throw new Error('Could not find the \"- Pedido (...)\" paragraph.');\n}\n\n// Search only inside that paragraph so we never touch the other\n// \"id_carrinho\" occurrences elsewhere in the document (e.g. the\n// \"Pedido[id_carrinho] => Carrinho[id]\" relationship line).\nconst matches = pedidoPara.search(\"id_carrinho\", { matchCase: true });\nmatches.load(\"items\");\nawait context.sync();\n\nif (matches.items.length === 0) {\n  throw new Error('Could not find \"id_carrinho\" inside the Pedido paragraph.');\n}\n\nconst idCarrinhoRange = matches.items[0];\n\n// Insert the new attribute right after \"id_carrinho\" (i.e. before the\n// closing \")\"), then make sure the newly inserted text is not italic \u2014\n// unlike \"id_carrinho\" (a foreign key, shown in italics), \"situacao\" is a\n// regular column, like \"meio_pagamento\".\nconst inserted = idCarrinhoRange.insertText(\", situacao\", Word.InsertLocation.after);\ninserted.font.italic = false;\n\nawait context.sync();\n", "ps1": "# Commit: \"adicionando atributo de situacao para pedido - futuro enum\"\n# The \"Pedido (...)\" bullet line needs a new \", situacao\" column appended\n# right after \"id_carrinho\" (and before the closing parenthesis), matching\n# the plain (non-italic) formatting used for the other regular attributes\n# in that line (only \"id_carrinho\" itself stays italic, marking it as a\n# foreign key).\n\n$d = $word.ActiveDocument\n\n# Locate the \"- Pedido (...)\" definition paragraph explicitly instead of\n# assuming a fixed index, so the script is resilient to minor reflow.\n$pedidoPara = $null\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text\n    if ($t -like \"*Pedido (*\" -and $t -like \"*id_carrinho*\") {\n        $pedidoPara = $p\n        break\n    }\n}\n\nif ($pedidoPara -eq $null) {\n    throw \"Could not find the '- Pedido (...)' paragraph.\"\n}\n\n# Search only inside that paragraph's range so we never touch the other\n# \"id_carrinho\" occurrences elsewhere in the document (e.g. the\n# \"Pedido[id_carrinho] => Carrinho[id]\" relationship line).\n$rng = $pedidoPara.Range\n$rng.Find.ClearFormatting()\n$found = $rng.Find.Execute(\"id_carrinho\")\nif (-not $found) {\n    throw \"Could not find 'id_carrinho' inside the Pedido paragraph.\"\n}\n\n# Collapse to the end of the match (right before the closing \")\") and\n# insert the new attribute there.\n$rng.Collapse(0)   # wdCollapseEnd\n$rng.InsertAfter(\", situacao\")\n\n# Make sure the newly inserted text is not italic -- unlike \"id_carrinho\"\n# (a foreign key, shown in italics), \"situacao\" is a regular column, like\n# \"meio_pagamento\". $rng tracks the just-inserted text after InsertAfter.\n$rng.Font.Italic = 0\n"}
